$d = $word.ActiveDocument

# Locate the unique target sentence fragment containing "seconds" (the one
# about users returning to the system after a short while), so we don't
# collide with the unrelated "minutes" that already appears earlier in the
# document.
$anchor = $d.Content.Duplicate
$anchor.Find.ClearFormatting()
$found = $anchor.Find.Execute(
    "seconds before beginning to use the system.",
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0
)

if (-not $found) {
    Write-Output "Target sentence not found"
} else {
    $secondsStart = $anchor.Start
    $secondsEnd = $secondsStart + "seconds".Length

    # --- Replace "seconds" with "minutes", splitting it into its own run by
    #     toggling a (reverted) character-format change around the edit. ---
    $wordRange = $d.Range($secondsStart, $secondsEnd)
    $wordRange.Bold = 1

    $replaceRange = $d.Range($secondsStart, $secondsEnd)
    $replaceRange.Text = "minutes"

    $newWordEnd = $secondsStart + "minutes".Length
    $unboldRange = $d.Range($secondsStart, $newWordEnd)
    $unboldRange.Bold = 0

    # --- Split the trailing phrase (" before beginning to use the
    #     system.") into its own run too, so it is no longer glued to the
    #     run that precedes "minutes". ---
    $tailPhrase = " before beginning to use the system."
    $tailStart = $newWordEnd
    $tailEnd = $tailStart + $tailPhrase.Length

    $tailRange = $d.Range($tailStart, $tailEnd)
    $tailRange.Bold = 1

    # Route the (otherwise unchanged) tail text through a dummy placeholder
    # so the engine registers a genuine text replacement here too - assigning
    # identical text back is treated as a no-op and would leave this text
    # glued to its neighboring run.
    $placeholder = "ZzPlaceholderTextzZ"
    $tailReplace = $d.Range($tailStart, $tailEnd)
    $tailReplace.Text = $placeholder

    $placeholderEnd = $tailStart + $placeholder.Length
    $tailRestore = $d.Range($tailStart, $placeholderEnd)
    $tailRestore.Text = $tailPhrase

    $tailUnbold = $d.Range($tailStart, $tailStart + $tailPhrase.Length)
    $tailUnbold.Bold = 0

    Write-Output "Replaced 'seconds' with 'minutes' in the target sentence"
}
